$wb = $excel.ActiveWorkbook

$wsRepay = $wb.Worksheets.Item("Repayment schedule")
$wsTrans = $wb.Worksheets.Item("Transactions")

# Insert a new blank column before column N ("Late") on the "Repayment schedule"
# sheet. This shifts the old N/O/P columns (Late / heading / Outstanding) one
# column to the right, becoming O/P/Q, and leaves the freshly inserted column N
# blank (inheriting the width/format of the column immediately to its left,
# i.e. the "In Advance" column).
$wsRepay.Columns("N:N").Insert()
$wsRepay.Columns("N:N").ColumnWidth = 9.83

# The "Transactions" sheet keeps its own saved selection, it's just no longer
# the active tab.
[void]$wsTrans.Select()
[void]$wsTrans.Range("D3").Select()

# Select the "Repayment schedule" sheet last (it becomes the tab shown when
# the workbook is opened) and move the selection there, as happened during
# the editing session.
[void]$wsRepay.Select()
[void]$wsRepay.Range("K13").Select()
